{"js": "// Update the worksheet date header and the 25 division problems in the\n// table to a new generated set (commit \"Update master to output generated\n// at c986bee\").\n\n// 1) Date header: first body paragraph, outside the table.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nif (paragraphs.items.length > 0 && paragraphs.items[0].text === \"2024-09-12 Thursday\") {\n  paragraphs.items[0].insertText(\"2024-09-13 Friday\", \"Replace\");\n}\n\n// 2) The 25 \"NN\u00f7N=\" problems live in a 5-column table; the five rows that\n// carry data are grid rows 0, 4, 8, 12, 16 (the rows in between are blank\n// spacer rows). Address each cell directly by (row, col) so the repeated\n// value \"99\u00f74=\" (both a source and a result elsewhere) can't collide with\n// a global text search/replace.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [gridRow, col, expectedOldValue, newValue]\nconst replacements = [\n  [0, 0, \"53\u00f78=\", \"99\u00f74=\"],\n  [0, 1, \"66\u00f79=\", \"54\u00f78=\"],\n  [0, 2, \"97\u00f77=\", \"40\u00f76=\"],\n  [0, 3, \"70\u00f73=\", \"59\u00f78=\"],\n  [0, 4, \"44\u00f73=\", \"85\u00f76=\"],\n  [4, 0, \"32\u00f73=\", \"91\u00f72=\"],\n  [4, 1, \"13\u00f77=\", \"74\u00f73=\"],\n  [4, 2, \"99\u00f74=\", \"30\u00f77=\"],\n  [4, 3, \"89\u00f77=\", \"95\u00f79=\"],\n  [4, 4, \"62\u00f73=\", \"93\u00f72=\"],\n  [8, 0, \"89\u00f74=\", \"53\u00f74=\"],\n  [8, 1, \"90\u00f74=\", \"69\u00f73=\"],\n  [8, 2, \"77\u00f76=\", \"84\u00f75=\"],\n  [8, 3, \"59\u00f74=\", \"95\u00f77=\"],\n  [8, 4, \"10\u00f75=\", \"68\u00f73=\"],\n  [12, 0, \"11\u00f74=\", \"30\u00f73=\"],\n  [12, 1, \"63\u00f79=\", \"56\u00f78=\"],\n  [12, 2, \"86\u00f79=\", \"43\u00f76=\"],\n  [12, 3, \"42\u00f74=\", \"77\u00f74=\"],\n  [12, 4, \"10\u00f78=\", \"36\u00f72=\"],\n  [16, 0, \"38\u00f75=\", \"48\u00f73=\"],\n  [16, 1, \"67\u00f76=\", \"68\u00f72=\"],\n  [16, 2, \"64\u00f72=\", \"51\u00f79=\"],\n  [16, 3, \"49\u00f77=\", \"91\u00f73=\"],\n  [16, 4, \"45\u00f73=\", \"11\u00f73=\"]\n];\n\nconst cells = replacements.map(([row, col]) => table.getCell(row, col));\ncells.forEach((cell) => cell.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [, , oldVal, newVal] = replacements[i];\n  const cell = cells[i];\n  if (cell.value === oldVal) {\n    cell.insertText(newVal, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date header and the 25 division problems in the\n# table to a new generated set (commit \"Update master to output generated\n# at c986bee\").\n\n$d = $word.ActiveDocument\n\n# 1) Date header: first body paragraph, outside the table.\n# Paragraph.Range.Text carries the trailing paragraph mark (\\r); compare\n# only the visible text.\n$dateParagraph = $d.Paragraphs.Item(1)\n$currentDate = $dateParagraph.Range.Text.TrimEnd([char]13)\nif ($currentDate -eq \"2024-09-12 Thursday\") {\n    $dateParagraph.Range.Text = \"2024-09-13 Friday\"\n}\n\n# 2) The 25 \"NN\u00f7N=\" problems live in a 5-column table; the five rows that\n# carry data are (1-based) table rows 1, 5, 9, 13, 17 (the rows in between\n# are blank spacer rows). Address each cell directly by (row, col) so the\n# repeated value \"99\u00f74=\" (both a source and a result elsewhere) can't\n# collide with a global find/replace.\n$table = $d.Tables.Item(1)\n\n# Each entry: row, col, expected old value, new value\n$replacements = @(\n    @(1, 1, \"53\u00f78=\", \"99\u00f74=\"),\n    @(1, 2, \"66\u00f79=\", \"54\u00f78=\"),\n    @(1, 3, \"97\u00f77=\", \"40\u00f76=\"),\n    @(1, 4, \"70\u00f73=\", \"59\u00f78=\"),\n    @(1, 5, \"44\u00f73=\", \"85\u00f76=\"),\n    @(5, 1, \"32\u00f73=\", \"91\u00f72=\"),\n    @(5, 2, \"13\u00f77=\", \"74\u00f73=\"),\n    @(5, 3, \"99\u00f74=\", \"30\u00f77=\"),\n    @(5, 4, \"89\u00f77=\", \"95\u00f79=\"),\n    @(5, 5, \"62\u00f73=\", \"93\u00f72=\"),\n    @(9, 1, \"89\u00f74=\", \"53\u00f74=\"),\n    @(9, 2, \"90\u00f74=\", \"69\u00f73=\"),\n    @(9, 3, \"77\u00f76=\", \"84\u00f75=\"),\n    @(9, 4, \"59\u00f74=\", \"95\u00f77=\"),\n    @(9, 5, \"10\u00f75=\", \"68\u00f73=\"),\n    @(13, 1, \"11\u00f74=\", \"30\u00f73=\"),\n    @(13, 2, \"63\u00f79=\", \"56\u00f78=\"),\n    @(13, 3, \"86\u00f79=\", \"43\u00f76=\"),\n    @(13, 4, \"42\u00f74=\", \"77\u00f74=\"),\n    @(13, 5, \"10\u00f78=\", \"36\u00f72=\"),\n    @(17, 1, \"38\u00f75=\", \"48\u00f73=\"),\n    @(17, 2, \"67\u00f76=\", \"68\u00f72=\"),\n    @(17, 3, \"64\u00f72=\", \"51\u00f79=\"),\n    @(17, 4, \"49\u00f77=\", \"91\u00f73=\"),\n    @(17, 5, \"45\u00f73=\", \"11\u00f73=\")\n)\n\nforeach ($entry in $replacements) {\n    $row = $entry[0]\n    $col = $entry[1]\n    $oldVal = $entry[2]\n    $newVal = $entry[3]\n    $cell = $table.Cell($row, $col)\n    $cellRange = $cell.Range\n    # Cell.Range.Text carries the trailing cell-mark (\\r\\a); compare only\n    # the visible text.\n    $current = $cellRange.Text.TrimEnd([char]13, [char]7)\n    if ($current -eq $oldVal) {\n        $cellRange.Text = $newVal\n    }\n}\n"}
